$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
$ws.Range("D2").Value = "27.508.77"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "1.844.72"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.01"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4599"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3860"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.04"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07918"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.000"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.51"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.860.33"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.973"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.142"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.37"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06699"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001035"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").Value = "27.524.18"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.399"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.93"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").Value = "2.074.04"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.02"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.52"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.117"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.431"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.44"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9737"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09393"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.620"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.292"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.335"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06002"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02224"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.287"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.184"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.008"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5903"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1862"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.33"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.243"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5578"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.11"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.909"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06697"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.20"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.047"
$ws.Range("E51").Value = "  -1.56%  "
